$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = 5373.000859812833
$ws.Range("E4").Value = 2904.218707803637
$ws.Range("F4").Value = 0.5405207971444107
$ws.Range("G4").Value = 1.850067574241423
$ws.Range("H4").Value = 130.4770519271406
$ws.Range("I4").Value = 15.62149409304038
$ws.Range("J4").Value = 10.94168601975434
$ws.Range("K4").Value = 259.5308018387295
$ws.Range("L4").Value = 251.6561699544545
$ws.Range("M4").Value = 126.7164286980405
$ws.Range("N4").Value = 2.784134542991524
$ws.Range("O4").Value = 4.935462251400168
$ws.Range("P4").Value = 0.007548424894594064
$ws.Range("Q4").Value = 348.7826346189249
$ws.Range("R4").Value = 522.286516635213
$ws.Range("S4").Value = 21.98354299389757
$ws.Range("T4").Value = -2606.75690399982
$ws.Range("U4").Value = -307.4944196094075
$ws.Range("V4").Value = -218.8261719701923
$ws.Range("W4").Value = -4841.833402155666
$ws.Range("X4").Value = -2012.042057325598
$ws.Range("Y4").Value = -5011.139856095193
$ws.Range("D5").Value = 5373.000859812833
$ws.Range("E5").Value = 3158.894986740313
$ws.Range("F5").Value = 0.5879200597876607
$ws.Range("G5").Value = 1.700911515693427
$ws.Range("H5").Value = 141.2996040211729
$ws.Range("I5").Value = 15.81421215635055
$ws.Range("J5").Value = 11.86395487017035
$ws.Range("K5").Value = 281.0072454321198
$ws.Range("L5").Value = 271.7145184960682
$ws.Range("M5").Value = 152.3988583693281
$ws.Range("N5").Value = 2.784134542991524
$ws.Range("O5").Value = 4.935462251400168
$ws.Range("P5").Value = 0.007548424894594064
$ws.Range("Q5").Value = 348.7826346189249
$ws.Range("R5").Value = 522.286516635213
$ws.Range("S5").Value = 21.98354299389757
$ws.Range("T5").Value = -2823.207945880466
$ws.Range("U5").Value = -311.3487808756108
$ws.Range("V5").Value = -237.2715489785123
$ws.Range("W5").Value = -5271.362274023471
$ws.Range("X5").Value = -2525.690650751349
$ws.Range("Y5").Value = -5412.306826927466
$ws.Range("D6").Value = 5373.000859812833
$ws.Range("E6").Value = 3413.440785474144
$ws.Range("F6").Value = 0.6352950380121566
$ws.Range("G6").Value = 1.574071793680316
$ws.Range("H6").Value = 152.1166113162471
$ws.Range("I6").Value = 16.00683148296594
$ws.Range("J6").Value = 12.78575120769437
$ws.Range("K6").Value = 302.472685838351
$ws.Range("L6").Value = 291.7625903938897
$ws.Range("M6").Value = 178.0681299697608
$ws.Range("N6").Value = 2.784134542991524
$ws.Range("O6").Value = 4.935462251400168
$ws.Range("P6").Value = 0.007548424894594064
$ws.Range("Q6").Value = 348.7826346189249
$ws.Range("R6").Value = 522.286516635213
$ws.Range("S6").Value = 21.98354299389757
$ws.Range("T6").Value = -3039.548091781951
$ws.Range("U6").Value = -315.2011674079185
$ws.Range("V6").Value = -255.7074757289929
$ws.Range("W6").Value = -5700.671082148096
$ws.Range("X6").Value = -3039.076082760002
$ws.Range("Y6").Value = -5813.268264883896
$ws.Range("D7").Value = 5373.000859812833
$ws.Range("E7").Value = 3667.856204257347
$ws.Range("F7").Value = 0.6826457504764135
$ws.Range("G7").Value = 1.464888632650401
$ws.Range("H7").Value = 162.9280780724803
$ws.Range("I7").Value = 16.19935214875841
$ws.Range("J7").Value = 13.70707539536306
$ws.Range("K7").Value = 323.9271315115038
$ws.Range("L7").Value = 311.8003935443703
$ws.Range("M7").Value = 203.7242536088452
$ws.Range("N7").Value = 2.784134542991524
$ws.Range("O7").Value = 4.935462251400168
$ws.Range("P7").Value = 0.007548424894594064
$ws.Range("Q7").Value = 348.7826346189249
$ws.Range("R7").Value = 522.286516635213
$ws.Range("S7").Value = 21.98354299389757
$ws.Range("T7").Value = -3255.777426906614
$ws.Range("U7").Value = -319.051580723768
$ws.Range("V7").Value = -274.1339594823667
$ws.Range("W7").Value = -6129.759995611152
$ws.Range("X7").Value = -3552.19855554169
$ws.Range("Y7").Value = -6214.024327893509
$ws.Range("D8").Value = 5373.000859812833
$ws.Range("E8").Value = 3922.141343234107
$ws.Range("F8").Value = 0.7299722158188401
$ws.Range("G8").Value = 1.369915153384651
$ws.Range("H8").Value = 173.7340085455799
$ws.Range("I8").Value = 16.39177422951252
$ws.Range("J8").Value = 14.62792779584288
$ws.Range("K8").Value = 345.3705908970442
$ws.Range("L8").Value = 331.8279358339496
$ws.Range("M8").Value = 229.367239384912
$ws.Range("N8").Value = 2.784134542991524
$ws.Range("O8").Value = 4.935462251400168
$ws.Range("P8").Value = 0.007548424894594064
$ws.Range("Q8").Value = 348.7826346189249
$ws.Range("R8").Value = 522.286516635213
$ws.Range("S8").Value = 21.98354299389757
$ws.Range("T8").Value = -3471.896036368606
$ws.Range("U8").Value = -322.9000223388502
$ws.Range("V8").Value = -292.551007491963
$ws.Range("W8").Value = -6558.629183321958
$ws.Range("X8").Value = -4065.058271063026
$ws.Range("Y8").Value = -6614.575173685094
